$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$normalStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = '31.155.32'
$ws.Range("E2").Value = '  +4.17%  '

$ws.Range("D3").Value = '1.701.22'
$ws.Range("E3").Value = '  +3.86%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.11'
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = '  +2.77%  '

$ws.Range("E6").Value = '  +2.67%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.67'
$ws.Range("D8").Style = $normalStyle
$ws.Range("E8").Value = '  +2.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.30'
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = '  +3.37%  '

$ws.Range("E10").Value = '  +3.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0642'
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = '  +5.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0912'
$ws.Range("D12").Style = $normalStyle
$ws.Range("E12").Value = '  +1.21%  '

$ws.Range("D13").Value = '1.950.22'
$ws.Range("E13").Value = '  +4.23%  '

$ws.Range("D14").Value = '1.716.37'
$ws.Range("E14").Value = '  +6.25%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.613'
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = '  +4.26%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.24'
$ws.Range("D16").Style = $normalStyle
$ws.Range("E16").Value = '  +8.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.15'
$ws.Range("D17").Style = $normalStyle
$ws.Range("E17").Value = '  +7.19%  '

$ws.Range("D18").Value = '31.158.08'
$ws.Range("E18").Value = '  +4.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.09'
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = '  +3.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '248.44'
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = '  +3.16%  '

$ws.Range("D21").Value = '0.0₃0722'
$ws.Range("E21").Value = '  +2.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.28'
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = '  +3.26%  '

$ws.Range("E24").Value = '  +2.12%  '

$ws.Range("E25").Value = '  -0.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.68'
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = '  +0.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.98'
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = '  +2.93%  '

$ws.Range("E28").Value = '  +2.76%  '

$ws.Range("E29").Value = '  +1.75%  '

$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.67'
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = '  +8.45%  '

$ws.Range("E32").Value = '  +2.76%  '

$ws.Range("E33").Value = '  +3.73%  '

$ws.Range("E34").Value = '  +5.65%  '

$ws.Range("D35").Value = '1.514.11'
$ws.Range("E35").Value = '  +6.27%  '

$ws.Range("E36").Value = '  +2.81%  '

$ws.Range("E37").Value = '  +1.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '83.30'
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = '  +8.97%  '

$ws.Range("E39").Value = '  +9.68%  '

$ws.Range("E40").Value = '  +4.58%  '

$ws.Range("E41").Value = '  -2.61%  '

$ws.Range("E42").Value = '  +0.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.04'
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = '  +3.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.850'
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = '  +1.91%  '

$ws.Range("E45").Value = '  +1.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.03'
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = '  +2.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("E48").Value = '  +6.73%  '

$ws.Range("E49").Value = '  +3.72%  '

$ws.Range("D50").Value = '1.839.39'
$ws.Range("E50").Value = '  +3.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '94.33'
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = '  +1.42%  '
